# Fix bug: deburring field change UI
# The "deburring" (and, on BTMI016, "deburringChanged") columns used to show
# the raw boolean-ish strings "true"/"false". The UI now shows descriptive
# text instead: "Deburring (one-sided)" for true, "No deburring" for false.
# A leading "'" forces PowerShell/COM to write the value as literal text
# (not a boolean/formula) and, importantly, keeps the existing cell style
# (which uses quotePrefix) instead of Excel minting a brand-new style.

$wb = $excel.ActiveWorkbook

$deburringOn  = "'Deburring (one-sided)"
$deburringOff = "'No deburring"

# ----------------------------------------------------------------------
# FPA011
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("FPA011")
$ws.Range("H2").Value = $deburringOn
$ws.Range("H3").Value = $deburringOff
$ws.Range("H4").Value = $deburringOn
$ws.Range("H5").Value = $deburringOff
$ws.Columns.Item(8).ColumnWidth = 21.1666666667
$ws.Range("H3").Select()

# ----------------------------------------------------------------------
# FPA012-013-015-017
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("FPA012-013-015-017")
$ws.Range("D2").Value = 50
$ws.Range("H2").Value = $deburringOn
$ws.Range("D3").Value = 60
$ws.Range("H3").Value = $deburringOff
$ws.Range("D4").Value = 40
$ws.Range("H4").Value = $deburringOn
$ws.Range("D5").Value = 30
$ws.Range("H5").Value = $deburringOff
$ws.Columns.Item(8).ColumnWidth = 21.1666666667
$ws.Range("H2").Select()

# ----------------------------------------------------------------------
# FPA014-016-020
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("FPA014-016-020")
$ws.Range("H2").Value = $deburringOn
$ws.Range("H3").Value = $deburringOff
$ws.Columns.Item(8).ColumnWidth = 21.1666666667
$ws.Range("H3").Select()

# ----------------------------------------------------------------------
# FPA018-019
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("FPA018-019")
$ws.Range("H2").Value = $deburringOn
$ws.Range("H3").Value = $deburringOff
$ws.Range("H4").Value = $deburringOn
$ws.Range("H5").Value = $deburringOff
$ws.Columns.Item(8).ColumnWidth = 20.1666666667
$ws.Range("H3").Select()
$excel.ActiveWindow.Zoom = 70

# ----------------------------------------------------------------------
# BTMI010
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BTMI010")
$ws.Range("H2").Value = $deburringOn
$ws.Range("H3").Value = $deburringOff
$ws.Range("H4").Value = $deburringOn
$ws.Range("H5").Value = $deburringOff
$ws.Columns.Item(8).ColumnWidth = 21.1666666667
$ws.Range("H3").Select()

# ----------------------------------------------------------------------
# BTMI016 (left active / selected last, matches the saved workbook state)
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BTMI016")
$ws.Range("H2").Value = $deburringOn
$ws.Range("T2").Value = $deburringOff
$ws.Range("H3").Value = $deburringOff
$ws.Range("T3").Value = $deburringOn
$ws.Range("H4").Value = $deburringOn
$ws.Range("T4").Value = $deburringOff
$ws.Range("H5").Value = $deburringOff
$ws.Range("T5").Value = $deburringOn
$ws.Range("U19").Select()
